$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.894.24"
$ws.Range("E2").Value = "  -1.20%  "
$ws.Range("D3").Value = "'2.310.95"
$ws.Range("E3").Value = "  +0.05%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'302.71"
$ws.Range("E5").Value = "  -1.92%  "
$ws.Range("D6").Value = "'100.14"
$ws.Range("E6").Value = "  -4.80%  "
$ws.Range("D7").Value = "'0.504"
$ws.Range("E7").Value = "  -4.13%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("E9").Value = "  -2.61%  "
$ws.Range("D10").Value = "'34.86"
$ws.Range("E10").Value = "  -2.76%  "
$ws.Range("D11").Value = "'0.0791"
$ws.Range("E11").Value = "  -2.37%  "
$ws.Range("E12").Value = "  +0.28%  "
$ws.Range("D13").Value = "'6.72"
$ws.Range("E13").Value = "  -3.37%  "
$ws.Range("D14").Value = "'2.670.28"
$ws.Range("E14").Value = "  -0.02%  "
$ws.Range("D15").Value = "'15.65"
$ws.Range("E15").Value = "  +3.61%  "
$ws.Range("D16").Value = "'2.294.77"
$ws.Range("E16").Value = "  -0.54%  "
$ws.Range("D17").Value = "'0.797"
$ws.Range("E17").Value = "  -0.50%  "
$ws.Range("D18").Value = "'42.819.61"
$ws.Range("E18").Value = "  -1.23%  "
$ws.Range("E19").Value = "  -1.17%  "
$ws.Range("D20").Value = "'0.0₃0905"
$ws.Range("E20").Value = "  -1.76%  "
$ws.Range("E21").Value = "  -2.33%  "
$ws.Range("D22").Value = "'67.92"
$ws.Range("D23").Value = "'235.77"
$ws.Range("E23").Value = "  -1.89%  "
$ws.Range("E24").Value = "  -2.72%  "
$ws.Range("E25").Value = "  -3.33%  "
$ws.Range("D27").Value = "'24.81"
$ws.Range("E27").Value = "  -0.94%  "
$ws.Range("E28").Value = "  -1.74%  "
$ws.Range("D29").Value = "'34.53"
$ws.Range("E29").Value = "  -4.59%  "
$ws.Range("D30").Value = "'165.37"
$ws.Range("E30").Value = "  +1.87%  "
$ws.Range("D31").Value = "'9.11"
$ws.Range("E31").Value = "  -4.86%  "
$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "  +0.09%  "
$ws.Range("D33").Value = "'5.01"
$ws.Range("E33").Value = "  -4.17%  "
$ws.Range("D34").Value = "'2.42"
$ws.Range("E34").Value = "  -4.82%  "
$ws.Range("D35").Value = "'4.48"
$ws.Range("E35").Value = "  -2.71%  "
$ws.Range("D36").Value = "'16.75"
$ws.Range("E36").Value = "  -8.45%  "
$ws.Range("D37").Value = "'0.0697"
$ws.Range("E37").Value = "  -4.70%  "
$ws.Range("E38").Value = "  -3.65%  "
$ws.Range("E40").Value = "  -4.63%  "
$ws.Range("E41").Value = "  -3.55%  "
$ws.Range("D42").Value = "'2.51"
$ws.Range("E42").Value = "  -0.16%  "
$ws.Range("D43").Value = "'1.975.98"
$ws.Range("E43").Value = "  +0.62%  "
$ws.Range("D44").Value = "'0.0280"
$ws.Range("E44").Value = "  -2.85%  "
$ws.Range("D45").Value = "'18.45"
$ws.Range("E45").Value = "  -1.35%  "
$ws.Range("D46").Value = "'10.26"
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("E47").Value = "  -5.44%  "
$ws.Range("E48").Value = "  -3.72%  "
$ws.Range("D49").Value = "'2.534.79"
$ws.Range("E49").Value = "  -0.04%  "
$ws.Range("E50").Value = "  -3.35%  "
$ws.Range("E51").Value = "  +0.29%  "
